# Update the "Förändrad" (Changed) date column (C) from 2023-10-03 (45202)
# to 2023-10-04 (45203) for all data rows (2 through 360).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C360").Value = 45203
